# Fix Training Data Issue
# The "Date" column (BF) held a malformed label like "2-17-2013-14"
# (month-day mashed together with the season string). Correct it to the
# real ISO game date "2014-02-17" for every data row (BF2:BF31), keeping
# the values as plain text (not reinterpreted as Excel date serials) and
# without disturbing the existing cell formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")

# Force text interpretation so Excel doesn't coerce "2014-02-17" into a
# date serial number when the value is assigned below.
$dateRange.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2014-02-17"
}

# Restore the default "Normal" style so these cells keep the same
# (unstyled) appearance they had before the text-number-format tweak.
$dateRange.Style = "Normal"
